# export latest mex files to xlsx. update segments and articles_db.
# Appends 5 new coded-segment rows (145-149) to Sheet1, mirroring the
# formatting of the last existing row (144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bullet = "$([char]0x25CF)"

# --- Set cell values first -------------------------------------------------
# Columns D and I use a "General" number format, so numeric-looking text
# (doc ids like "20339", or "2013") would otherwise be auto-coerced into a
# number by Excel. A leading apostrophe forces these to be stored as text,
# matching the source data (doc ids and years are text in this sheet).

# Row 145
$ws.Range("A145").Value = $bullet
$ws.Range("B145").Value = "'"
$ws.Range("C145").Value = "'"
$ws.Range("D145").Value = "'20339"
$ws.Range("E145").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F145").Value = "1: 3713"
$ws.Range("G145").Value = "1: 3733"
$ws.Range("H145").Value = 0
$ws.Range("I145").Value = "Clostridium difficile"
$ws.Range("J145").Value = 21
$ws.Range("K145").Value = 0.093993
$ws.Range("L145").Value = "Sonia"
$ws.Range("M145").Value = "11/8/18 14:04:00"

# Row 146
$ws.Range("A146").Value = $bullet
$ws.Range("B146").Value = "'"
$ws.Range("C146").Value = "'"
$ws.Range("D146").Value = "'15902"
$ws.Range("E146").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F146").Value = "1: 1647"
$ws.Range("G146").Value = "1: 1668"
$ws.Range("H146").Value = 0
$ws.Range("I146").Value = "Nocardia transvalensis"
$ws.Range("J146").Value = 22
$ws.Range("K146").Value = 0.149976
$ws.Range("L146").Value = "Sonia"
$ws.Range("M146").Value = "11/8/18 14:05:00"

# Row 147
$ws.Range("A147").Value = $bullet
$ws.Range("B147").Value = "'"
$ws.Range("C147").Value = "'"
$ws.Range("D147").Value = "'18631"
$ws.Range("E147").Value = "Event year"
$ws.Range("F147").Value = "1: 1332"
$ws.Range("G147").Value = "1: 1336"
$ws.Range("H147").Value = 0
$ws.Range("I147").Value = "'2007."
$ws.Range("J147").Value = 5
$ws.Range("K147").Value = 0.029303
$ws.Range("L147").Value = "Sonia"
$ws.Range("M147").Value = "11/12/18 12:31:00"

# Row 148
$ws.Range("A148").Value = $bullet
$ws.Range("B148").Value = "'"
$ws.Range("C148").Value = "'"
$ws.Range("D148").Value = "'22167"
$ws.Range("E148").Value = "Event month"
$ws.Range("F148").Value = "3: 1441"
$ws.Range("G148").Value = "3: 1447"
$ws.Range("H148").Value = 0
$ws.Range("I148").Value = "October"
$ws.Range("J148").Value = 7
$ws.Range("K148").Value = 0.018916
$ws.Range("L148").Value = "Sonia"
$ws.Range("M148").Value = "11/12/18 12:32:00"

# Row 149
$ws.Range("A149").Value = $bullet
$ws.Range("B149").Value = "'"
$ws.Range("C149").Value = "'"
$ws.Range("D149").Value = "'22167"
$ws.Range("E149").Value = "Event year"
$ws.Range("F149").Value = "3: 1449"
$ws.Range("G149").Value = "3: 1452"
$ws.Range("H149").Value = 0
$ws.Range("I149").Value = "'2013"
$ws.Range("J149").Value = 4
$ws.Range("K149").Value = 0.010809
$ws.Range("L149").Value = "Sonia"
$ws.Range("M149").Value = "11/12/18 12:32:00"

# --- Apply formatting -------------------------------------------------------
# Copy the formatting (fill/border/font/number-format) from the last
# pre-existing row down onto the 5 new rows so the new cells match the
# established style, without disturbing the values/types set above.
$ws.Range("A144:M144").Copy()
$ws.Range("A145:M149").PasteSpecial(-4122)

# Match the row height used by the other data rows.
$templateHeight = $ws.Rows.Item(144).RowHeight
$ws.Rows.Item(145).RowHeight = $templateHeight
$ws.Rows.Item(146).RowHeight = $templateHeight
$ws.Rows.Item(147).RowHeight = $templateHeight
$ws.Rows.Item(148).RowHeight = $templateHeight
$ws.Rows.Item(149).RowHeight = $templateHeight
